$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume figures (rows 2-51) plus a row-48/49
# swap (BinanceUSD now ranks above SynthetixNetwork).

$ws.Range("D2").Value = "41.785.04"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.222.07"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.16%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.86%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0969"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.38%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "2.554.05"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.889"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "2.207.23"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "41.813.40"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "0.0₃0967"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.03%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.45%  "
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0309"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.203"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("E50").Value = "  +5.91%  "
$ws.Range("E51").Value = "  +5.29%  "
